$d = $word.ActiveDocument

# Replace the text "Primeiro teste commit" with "Testando git."
$d.Content.Find.Execute("Primeiro teste commit", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Testando git.", 2)
